# Insert a new first data column (FY2018, period-ending 2018-12-31) into the
# "MD" yearly financials sheet, pushing the existing FY2017..FY2011 columns
# (D..K) one column to the right (E..L) and filling in the new column D with
# the latest year's figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank column before column D; existing D:K data shifts to E:L.
$ws.Columns("D").Insert()

# Copy number formats/styles from the (now shifted) neighboring column E into
# the new column D so the new column renders the same way (date row, plain
# numbers, etc.) as the rest of the table. Done per contiguous data block so
# we don't materialize phantom cells on blank separator rows (6, 36, 37, 78,
# 79) that have no cells at all in columns D:L.
$ws.Range("E7:E35").Copy()
$ws.Range("D7:D35").PasteSpecial(-4122)

$ws.Range("E38:E77").Copy()
$ws.Range("D38:D77").PasteSpecial(-4122)

$ws.Range("E80:E102").Copy()
$ws.Range("D80:D102").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# Period-ending header row (date serials for FY2018-12-31).
$ws.Range("D7").Value = 43465
$ws.Range("D38").Value = 43465
$ws.Range("D80").Value = 43465

# Income statement (rows 8-35)
$ws.Range("D8").Value = 3647100
$ws.Range("D9").Value = "NA"
$ws.Range("D10").Value = "NA"
$ws.Range("D12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("D15").Value = 111300
$ws.Range("D17").Value = 3201300
$ws.Range("D18").Value = 445800
$ws.Range("D20").Value = 11800
$ws.Range("D21").Value = 568900
$ws.Range("D22").Value = 88800
$ws.Range("D23").Value = 368800
$ws.Range("D24").Value = 100200
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = 268600
$ws.Range("D27").Value = 268600
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = "NA"
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = -11800
$ws.Range("D33").Value = 268600
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = 268600

# Balance sheet (rows 41-77)
$ws.Range("D41").Value = 36700
$ws.Range("D42").Value = 21900
$ws.Range("D43").Value = 542300
$ws.Range("D44").Value = 0
$ws.Range("D45").Value = 56500
$ws.Range("D46").Value = 657400
$ws.Range("D47").Value = 69700
$ws.Range("D48").Value = 133000
$ws.Range("D49").Value = 4971300
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 103500
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 5934900
$ws.Range("D57").Value = 31100
$ws.Range("D58").Value = 300
$ws.Range("D59").Value = 468900
$ws.Range("D60").Value = 500200
$ws.Range("D61").Value = 1974300
$ws.Range("D62").Value = 372600
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 2847000
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = 2094400
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 3087900
$ws.Range("D77").Value = 0

# Cash flow statement (rows 81-102)
$ws.Range("D81").Value = 268600
$ws.Range("D83").Value = 111300
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = 289900
$ws.Range("D91").Value = -48900
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = -142800
$ws.Range("D96").Value = 0
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = -170600
$ws.Range("D101").Value = 0
$ws.Range("D102").Value = -23500

# Keep the new column's width in line with the rest of the (best-fit) data
# columns it now sits alongside.
$ws.Columns("D").ColumnWidth = $ws.Columns("E").ColumnWidth
